$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 16740.666
$ws.Range("I62").Value = 24111
$ws.Range("J62").Value = 2000
$ws.Range("K62").Value = 24111
$ws.Range("L62").Value = 2000
$ws.Range("M62").Value = -23487
$ws.Range("N62").Value = -3248
$ws.Range("H65").Value = 16740.666
$ws.Range("I65").Value = 24111
$ws.Range("J65").Value = 2000
$ws.Range("K65").Value = 120555
$ws.Range("L65").Value = 10000
$ws.Range("M65").Value = -117435
$ws.Range("N65").Value = -16240
$ws.Range("H132").Value = 1541.25
$ws.Range("I132").Value = 1539.5714
$ws.Range("K132").Value = 4618.7142
$ws.Range("M132").Value = -2088.7142
$ws.Range("H138").Value = 1714.86
$ws.Range("J138").Value = 1930.2113
$ws.Range("L138").Value = 5790.6339
$ws.Range("N138").Value = -16070.6339
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 241895.83
$ws.Range("I2").Value = 308933
$ws.Range("K2").Value = 308933
$ws.Range("M2").Value = -308820
$ws.Range("H32").Value = 4032.842
$ws.Range("I32").Value = 2484.3828
$ws.Range("J32").Value = 12991.786
$ws.Range("K32").Value = 2484.3828
$ws.Range("L32").Value = 12991.786
$ws.Range("M32").Value = -2197.3828
$ws.Range("N32").Value = -13565.786
$ws.Range("H61").Value = 60192.785
$ws.Range("I61").Value = 67808.336
$ws.Range("J61").Value = 14499.5
$ws.Range("K61").Value = 67808.336
$ws.Range("L61").Value = 14499.5
$ws.Range("M61").Value = -67596.336
$ws.Range("N61").Value = -14923.5
$ws.Range("H74").Value = 1183.25
$ws.Range("I74").Value = 642
$ws.Range("J74").Value = 2807
$ws.Range("K74").Value = 642
$ws.Range("L74").Value = 2807
$ws.Range("M74").Value = 232
$ws.Range("N74").Value = -4555
$ws.Range("H77").Value = 1183.25
$ws.Range("I77").Value = 642
$ws.Range("J77").Value = 2807
$ws.Range("K77").Value = 3210
$ws.Range("L77").Value = 14035
$ws.Range("M77").Value = 1158
$ws.Range("N77").Value = -22771
$ws.Range("H97").Value = 1937.8948
$ws.Range("I97").Value = 2001.0769
$ws.Range("K97").Value = 2001.0769
$ws.Range("M97").Value = -1505.0769
$ws.Range("H116").Value = 241895.83
$ws.Range("I116").Value = 308933
$ws.Range("K116").Value = 308933
$ws.Range("M116").Value = -306639
$ws.Range("H132").Value = 1893.6774
$ws.Range("I132").Value = 1554.28
$ws.Range("J132").Value = 3307.8333
$ws.Range("K132").Value = 4662.84
$ws.Range("L132").Value = 9923.499899999999
$ws.Range("M132").Value = -2132.84
$ws.Range("N132").Value = -14983.4999
$ws.Range("H136").Value = 60192.785
$ws.Range("I136").Value = 67808.336
$ws.Range("J136").Value = 14499.5
$ws.Range("K136").Value = 203425.008
$ws.Range("L136").Value = 43498.5
$ws.Range("M136").Value = -200875.008
$ws.Range("N136").Value = -48598.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 241895.83
$ws.Range("I3").Value = 308933
$ws.Range("K3").Value = 308933
$ws.Range("M3").Value = -308819
$ws.Range("H22").Value = 766.6667
$ws.Range("I22").Value = 100
$ws.Range("J22").Value = 1100
$ws.Range("K22").Value = 100
$ws.Range("L22").Value = 1100
$ws.Range("M22").Value = 73
$ws.Range("N22").Value = -1446
$ws.Range("H86").Value = 286606.28
$ws.Range("I86").Value = 970.3333
$ws.Range("K86").Value = 970.3333
$ws.Range("M86").Value = 152.6667
$ws.Range("H89").Value = 286606.28
$ws.Range("I89").Value = 970.3333
$ws.Range("K89").Value = 4851.6665
$ws.Range("M89").Value = 764.3334999999997
$ws.Range("H105").Value = 2411.7585
$ws.Range("I105").Value = 2287.3635
$ws.Range("J105").Value = 2802.7144
$ws.Range("K105").Value = 2287.3635
$ws.Range("L105").Value = 2802.7144
$ws.Range("M105").Value = -540.3634999999999
$ws.Range("N105").Value = -6296.7144
$ws.Range("H134").Value = 8581.105
$ws.Range("I134").Value = 12146.417
$ws.Range("J134").Value = 2469.1428
$ws.Range("K134").Value = 36439.251
$ws.Range("L134").Value = 7407.428400000001
$ws.Range("M134").Value = -33904.251
$ws.Range("N134").Value = -12477.4284
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1488.3077
$ws.Range("J22").Value = 1739.7
$ws.Range("L22").Value = 1739.7
$ws.Range("N22").Value = -2439.7
$ws.Range("H31").Value = 2527.3076
$ws.Range("J31").Value = 3482.1667
$ws.Range("L31").Value = 3482.1667
$ws.Range("N31").Value = -4072.1667
$ws.Range("H34").Value = 2527.3076
$ws.Range("J34").Value = 3482.1667
$ws.Range("L34").Value = 3482.1667
$ws.Range("N34").Value = -3886.1667
$ws.Range("H58").Value = 2290216.2
$ws.Range("I58").Value = 3107435.2
$ws.Range("J58").Value = 2002.4
$ws.Range("K58").Value = 3107435.2
$ws.Range("L58").Value = 2002.4
$ws.Range("M58").Value = -3107232.2
$ws.Range("N58").Value = -2408.4
$ws.Range("H132").Value = 2159.6667
$ws.Range("I132").Value = 1348.4286
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 4045.2858
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -1515.2858
$ws.Range("N132").Value = -20057
$ws.Range("H134").Value = 2517.2727
$ws.Range("I134").Value = 2458.077
$ws.Range("K134").Value = 7374.231000000001
$ws.Range("M134").Value = -4839.231000000001
$ws.Range("H136").Value = 2290216.2
$ws.Range("I136").Value = 3107435.2
$ws.Range("J136").Value = 2002.4
$ws.Range("K136").Value = 9322305.600000001
$ws.Range("L136").Value = 6007.200000000001
$ws.Range("M136").Value = -9319755.600000001
$ws.Range("N136").Value = -11107.2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 4221.8887
$ws.Range("I110").Value = 1500
$ws.Range("K110").Value = 4500
$ws.Range("M110").Value = -410
$ws.Range("H113").Value = 33142.824
$ws.Range("J113").Value = 841.62964
$ws.Range("L113").Value = 2524.88892
$ws.Range("N113").Value = -6864.888919999999
$ws.Range("H131").Value = 22755.621
$ws.Range("J131").Value = 23376.611
$ws.Range("L131").Value = 70129.833
$ws.Range("N131").Value = -80209.833
$ws.Range("H136").Value = 1342.4706
$ws.Range("I136").Value = 1281.4667
$ws.Range("K136").Value = 3844.4001
$ws.Range("M136").Value = 1255.5999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 929.6177
$ws.Range("I97").Value = 941.1923
$ws.Range("K97").Value = 941.1923
$ws.Range("M97").Value = -445.1923
$ws.Range("H132").Value = 1376504
$ws.Range("I132").Value = 1750708.5
$ws.Range("K132").Value = 5252125.5
$ws.Range("M132").Value = -5249595.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5332.6924
$ws.Range("I16").Value = 5739.5835
$ws.Range("K16").Value = 5739.5835
$ws.Range("M16").Value = -5569.5835
$ws.Range("H136").Value = 3179.3076
$ws.Range("I136").Value = 2782.1
$ws.Range("J136").Value = 4503.3335
$ws.Range("K136").Value = 8346.299999999999
$ws.Range("L136").Value = 13510.0005
$ws.Range("M136").Value = -5796.299999999999
$ws.Range("N136").Value = -18610.0005
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 12385.571
$ws.Range("J96").Value = 16959.8
$ws.Range("L96").Value = 16959.8
$ws.Range("N96").Value = -19705.8
$ws.Range("H122").Value = 31070.926
$ws.Range("I122").Value = 41250.25
$ws.Range("K122").Value = 123750.75
$ws.Range("M122").Value = -121300.75
$ws.Range("H126").Value = 2052.2273
$ws.Range("I126").Value = 1832.5
$ws.Range("K126").Value = 5497.5
$ws.Range("M126").Value = -3027.5
$ws.Range("H130").Value = 33728.69
$ws.Range("J130").Value = 33728.69
$ws.Range("L130").Value = 33728.69
$ws.Range("N130").Value = -43768.69
$ws.Range("H132").Value = 1338.2
$ws.Range("I132").Value = 1137.1471
$ws.Range("J132").Value = 2477.5
$ws.Range("K132").Value = 3411.4413
$ws.Range("L132").Value = 7432.5
$ws.Range("M132").Value = -881.4412999999995
$ws.Range("N132").Value = -12492.5
$ws.Range("H136").Value = 32681218
$ws.Range("I136").Value = 50505630
$ws.Range("J136").Value = 3122.5
$ws.Range("K136").Value = 151516890
$ws.Range("L136").Value = 9367.5
$ws.Range("M136").Value = -151514340
$ws.Range("N136").Value = -14467.5
